# Insert a new weekly data row at row 651 (shifting the existing rows
# 651-666 down to 652-667) and populate it with the new Choclo record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 651; this shifts rows
# 651:666 down to 652:667 and extends the sheet dimension to R667.
$ws.Rows("651:651").Insert()

# Populate the newly inserted row 651 with the new record's values.
$ws.Range("A651").Value = 10
$ws.Range("B651").Value = "Vega Modelo de Temuco"
$ws.Range("C651").Value = "La Araucanía"
$ws.Range("D651").Value = 45239
$ws.Range("E651").Value = 9
$ws.Range("F651").Value = 100112024
$ws.Range("G651").Value = "Choclo"
$ws.Range("H651").Value = "Dulce o Americano"
$ws.Range("I651").Value = "Primera"
$ws.Range("J651").Value = 170
$ws.Range("K651").Value = 30000
$ws.Range("L651").Value = 35000
$ws.Range("M651").Value = 32647
$ws.Range("N651").Value = '$/malla 50 unidades'
$ws.Range("O651").Value = "Región de Arica y Parinacota"
$ws.Range("P651").Value = 653
$ws.Range("Q651").Value = 50
$ws.Range("R651").Value = "Hortaliza"
